$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date number-format from the last existing date cell (G64) down onto the
# new date cells (G65:G68) so they reuse the existing style record instead of
# creating a duplicate number format in styles.xml
$ws.Range("G64").Copy()
$ws.Range("G65:G68").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---- Row 65: survey_round 53, panel F, wave 17, received 2021-04-06 ----
$ws.Range("A65").Value = 3
$ws.Range("B65").Value = 0
$ws.Range("C65").Value = "uk"
$ws.Range("D65").Value = 53
$ws.Range("E65").Value = "F"
$ws.Range("F65").Value = 17
$ws.Range("G65").Value = "2021-04-06"
$ws.Range("H65").Value = "20-100590_PFW17_Final_ICUO"
$ws.Range("I65").Formula = "=C65&""_""&""sr""&TEXT(D65,""00"")&""_""&YEAR(G65)&TEXT(G65,""MM"")&TEXT(G65,""DD"")&""_p""&E65&""_wv""&TEXT(F65,""00"")&"""""
$ws.Range("J65").Value = 1

# ---- Row 66: survey_round 54, panel E, wave 18, received 2021-04-09 ----
$ws.Range("A66").Value = 3
$ws.Range("B66").Value = 0
$ws.Range("C66").Value = "uk"
$ws.Range("D66").Value = 54
$ws.Range("E66").Value = "E"
$ws.Range("F66").Value = 18
$ws.Range("G66").Value = "2021-04-09"
$ws.Range("H66").Value = "20-100562_PEW18_Final_ICUO"
$ws.Range("I66").Formula = "=C66&""_""&""sr""&TEXT(D66,""00"")&""_""&YEAR(G66)&TEXT(G66,""MM"")&TEXT(G66,""DD"")&""_p""&E66&""_wv""&TEXT(F66,""00"")&"""""
$ws.Range("J66").Value = 1

# ---- Row 67: survey_round 55, panel F, wave 18, received 2021-04-15 ----
$ws.Range("A67").Value = 3
$ws.Range("B67").Value = 0
$ws.Range("C67").Value = "uk"
$ws.Range("D67").Value = 55
$ws.Range("E67").Value = "F"
$ws.Range("F67").Value = 18
$ws.Range("G67").Value = "2021-04-15"
$ws.Range("H67").Value = "20-100590_PFW18_Final_ICUO"
$ws.Range("I67").Formula = "=C67&""_""&""sr""&TEXT(D67,""00"")&""_""&YEAR(G67)&TEXT(G67,""MM"")&TEXT(G67,""DD"")&""_p""&E67&""_wv""&TEXT(F67,""00"")&"""""
$ws.Range("J67").Value = 1

# ---- Row 68: survey_round 56, panel E, wave 19, received 2021-04-22 ----
$ws.Range("A68").Value = 3
$ws.Range("B68").Value = 0
$ws.Range("C68").Value = "uk"
$ws.Range("D68").Value = 56
$ws.Range("E68").Value = "E"
$ws.Range("F68").Value = 19
$ws.Range("G68").Value = "2021-04-22"
$ws.Range("H68").Value = "20-100562_PEW19_Final_ICUO"
$ws.Range("I68").Formula = "=C68&""_""&""sr""&TEXT(D68,""00"")&""_""&YEAR(G68)&TEXT(G68,""MM"")&TEXT(G68,""DD"")&""_p""&E68&""_wv""&TEXT(F68,""00"")&"""""
$ws.Range("J68").Value = 1

# Match the author's final selection (the new last row's received_final cell)
$ws.Range("J68").Select()
